# Player.xlsx / "Packet" sheet
# GameEnter response addition (PlayerPacket): replace the old
# "PlacedKingdomItemList : LIST:PlacedKingdomItemPacket" row (row 30) with
# "KingdomMap : KingdomMapPacket", and give every LIST/packet-typed field
# (rows 24-30) a default Value of "new()" in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the row 30 field from PlacedKingdomItemList/LIST:PlacedKingdomItemPacket
# to KingdomMap/KingdomMapPacket.
$ws.Range("A30").Value = "KingdomMap"
$ws.Range("B30").Value = "KingdomMapPacket"

# Add the default-value column (C) for the list/packet fields, rows 24-30.
$ws.Range("C24").Value = "new()"
$ws.Range("C25").Value = "new()"
$ws.Range("C26").Value = "new()"
$ws.Range("C27").Value = "new()"
$ws.Range("C28").Value = "new()"
$ws.Range("C29").Value = "new()"
$ws.Range("C30").Value = "new()"

# Leave the cursor where the author left it.
$ws.Range("C30").Select()
